# Add season-record columns (Wins / Losses / Ties) to the NYY_2000 sheet.
#
# The sheet currently ends at column AC (Unnamed: 28). We append three new
# columns: AD = Wins, AE = Losses, AF = Ties, with header labels in row 1
# (styled like the rest of the header row) and the team's 2000 season
# record (87 wins, 74 losses, 0 ties) repeated down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 48 }

# --- Header row ------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used for the other header cells (bold
# font, thin border, centered/top aligned) by copying the format from
# the preceding header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows ---------------------------------------------------------
# Every row in this sheet represents a New York Yankees player from the
# 2000 season, so every row gets the same team record.
$wins = 87
$losses = 74
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
